# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header / summary values
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 1256346
$ws.Range("C13").Value = 21
$ws.Range("F13").Value = 1

# ---------------------------------------------------------------------
# 2) Grow the workers table from 4 rows (16-19) to 21 rows (16-36).
#    Insert 17 blank rows at row 19 - this pushes the former row 19
#    (which carries the "last row" style) down to row 36, and the
#    footer block (formerly rows 24-25) down to rows 41-42.
# ---------------------------------------------------------------------
$ws.Rows.Item(19).Resize(17).Insert()

# Give the 17 freshly-inserted rows (19-35) the same formatting as the
# existing "middle" data row (row 18).
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Fill in the 21 rows of worker data (columns B:G); H:J stay blank.
# ---------------------------------------------------------------------
$data = @(
    @("45369441", "DORIS MARIA FUENTES ZUÑIGA", "2508", 56940, 1423500),
    @("45368418", "ANA VICTORIA DELGADO ROCHA", "2508", 63674, 1591848),
    @("45368149", "MILADIS CASTRO MARTINEZ", "2508", 63674, 1591848),
    @("45368997", "LINET VEGA MENDOZA", "2508", 63674, 1591848),
    @("23126816", "BERLIDES DEL CARMEN RUIZ CORTES", "2508", 56940, 1423500),
    @("22968820", "EDILMA ARNEDO TORRES", "2508", 56940, 1423500),
    @("45366594", "YADIRA ROCHA MARTINEZ", "2508", 56940, 1423500),
    @("45366621", "MARIA GUILLERMINA ROCHA BALBUENA", "2508", 63674, 1591848),
    @("45366745", "YORSELIS MALDONADO HERNANDEZ", "2508", 63674, 1591848),
    @("30655222", "NELLY DEL CARMEN MARTINEZ ANAYA", "2508", 63674, 1591848),
    @("45367257", "LIBIA ROSA DE ARCO TARRA", "2508", 63674, 1591848),
    @("45368061", "VERNUIL MIRANDA SOTO", "2508", 63674, 1591848),
    @("51622849", "ETEL CLEOPATRA TAPIA DE ARCO", "2508", 63674, 1591848),
    @("45371729", "BELINA FUENTES CASSIANI", "2508", 56940, 1423500),
    @("9039688", "JUAN REYES BLANCO", "2508", 56940, 1423500),
    @("1049942955", "YIPSY ALEJANDRA ALVAREZ FUENTES", "2508", 56940, 1423500),
    @("1049928787", "ROSARIO MERCEDES TORRES BOLIVAR", "2508", 56940, 1423500),
    @("1048943092", "YURANIS PAOLA CONTRERAS VALDES", "2508", 56940, 1423500),
    @("45370363", "YARIS IRIARTE ZAMBRANO", "2508", 56940, 1423500),
    @("45367827", "FRANCIA ELENA PEREZ PEREZ", "2508", 56940, 1423500),
    @("1049935694", "YESSICA PAOLA CARO MUNARRIS", "2508", 56940, 1423500)
)

$r = 16
foreach ($row in $data) {
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = $row[0]
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("E$r").Value = $row[2]
    $ws.Range("F$r").Value = $row[3]
    $ws.Range("G$r").Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 4) Column D got a bit wider to fit the longer names.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 36.6328125
